{"js": "// Update the two-digit multiplication problems in the practice table.\n// The document contains a table of \"AxB=\" prompts; this diff swaps each\n// prompt's text for a new pair of factors, in document order. We look up\n// each old prompt by its exact text and replace it with the new prompt.\n// Replacements are applied strictly in document order with a single\n// search+replace per pair, so an old value that happens to equal a later\n// pair's new value is never touched twice.\nconst replacements = [\n  [\"85\u00d790=\", \"35\u00d722=\"],\n  [\"69\u00d785=\", \"43\u00d731=\"],\n  [\"96\u00d718=\", \"21\u00d761=\"],\n  [\"90\u00d721=\", \"24\u00d737=\"],\n  [\"71\u00d766=\", \"69\u00d767=\"],\n  [\"21\u00d781=\", \"92\u00d713=\"],\n  [\"12\u00d775=\", \"30\u00d731=\"],\n  [\"16\u00d755=\", \"83\u00d784=\"],\n  [\"56\u00d764=\", \"70\u00d775=\"],\n  [\"43\u00d782=\", \"43\u00d723=\"],\n  [\"20\u00d782=\", \"93\u00d754=\"],\n  [\"36\u00d740=\", \"47\u00d757=\"],\n  [\"29\u00d781=\", \"93\u00d750=\"],\n  [\"62\u00d795=\", \"23\u00d761=\"],\n  [\"13\u00d764=\", \"15\u00d797=\"],\n  [\"59\u00d793=\", \"41\u00d745=\"],\n  [\"53\u00d798=\", \"31\u00d712=\"],\n  [\"54\u00d765=\", \"20\u00d782=\"],\n  [\"14\u00d781=\", \"76\u00d738=\"],\n  [\"85\u00d797=\", \"35\u00d726=\"],\n  [\"67\u00d779=\", \"81\u00d743=\"],\n  [\"28\u00d787=\", \"98\u00d799=\"],\n  [\"16\u00d792=\", \"48\u00d767=\"],\n  [\"28\u00d777=\", \"78\u00d743=\"],\n  [\"90\u00d798=\", \"65\u00d731=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Could not find text to replace: \" + oldText);\n  }\n\n  // Replace only the first (and expected only) occurrence so an identical\n  // old value appearing twice can't be double-consumed.\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Update the two-digit multiplication problems in the practice table.\n# The document contains a table of \"AxB=\" prompts; this diff swaps each\n# prompt's text for a new pair of factors, in document order. We look up\n# each old prompt by its exact text (Find.Execute) and replace it with the\n# new prompt text. Replacements are applied strictly in document order,\n# one exact match at a time (wdReplaceOne), so an old value that happens\n# to equal a later pair's new value is never touched twice.\n\n$olds = @(\"85\u00d790=\", \"69\u00d785=\", \"96\u00d718=\", \"90\u00d721=\", \"71\u00d766=\", \"21\u00d781=\", \"12\u00d775=\", \"16\u00d755=\", \"56\u00d764=\", \"43\u00d782=\", \"20\u00d782=\", \"36\u00d740=\", \"29\u00d781=\", \"62\u00d795=\", \"13\u00d764=\", \"59\u00d793=\", \"53\u00d798=\", \"54\u00d765=\", \"14\u00d781=\", \"85\u00d797=\", \"67\u00d779=\", \"28\u00d787=\", \"16\u00d792=\", \"28\u00d777=\", \"90\u00d798=\")\n\n$news = @(\"35\u00d722=\", \"43\u00d731=\", \"21\u00d761=\", \"24\u00d737=\", \"69\u00d767=\", \"92\u00d713=\", \"30\u00d731=\", \"83\u00d784=\", \"70\u00d775=\", \"43\u00d723=\", \"93\u00d754=\", \"47\u00d757=\", \"93\u00d750=\", \"23\u00d761=\", \"15\u00d797=\", \"41\u00d745=\", \"31\u00d712=\", \"20\u00d782=\", \"76\u00d738=\", \"35\u00d726=\", \"81\u00d743=\", \"98\u00d799=\", \"48\u00d767=\", \"78\u00d743=\", \"65\u00d731=\")\n\n$d = $word.ActiveDocument\n\nfor ($i = 0; $i -lt $olds.Count; $i++) {\n    $oldText = $olds[$i]\n    $newText = $news[$i]\n\n    $findRange = $d.Content\n    $found = $findRange.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 1)\n\n    if (-not $found) {\n        Write-Output \"WARNING: text not found: $oldText\"\n    }\n}\n"}
